$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 815.8570999999999
$ws.Range("I28").Value = 814.8333
$ws.Range("K28").Value = 814.8333
$ws.Range("M28").Value = -329.8333
$ws.Range("H33").Value = 142.07143
$ws.Range("I33").Value = 149
$ws.Range("J33").Value = 116.666664
$ws.Range("K33").Value = 149
$ws.Range("L33").Value = 116.666664
$ws.Range("M33").Value = 80
$ws.Range("N33").Value = -574.666664
$ws.Range("H34").Value = 10397.4
$ws.Range("I34").Value = 10496.75
$ws.Range("K34").Value = 10496.75
$ws.Range("M34").Value = -10293.75
$ws.Range("H36").Value = 10397.4
$ws.Range("I36").Value = 10496.75
$ws.Range("K36").Value = 10496.75
$ws.Range("M36").Value = -9781.75
$ws.Range("H53").Value = 372.2857
$ws.Range("I53").Value = 155
$ws.Range("K53").Value = 155
$ws.Range("M53").Value = 482
$ws.Range("H76").Value = 1500
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 1500
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 1500
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -2130
$ws.Range("H79").Value = 1500
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 1500
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 1500
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -3684
$ws.Range("H107").Value = 326.2
$ws.Range("I107").Value = 213.8
$ws.Range("K107").Value = 213.8
$ws.Range("M107").Value = 1706.2
$ws.Range("H112").Value = 4884.3335
$ws.Range("J112").Value = 4884.3335
$ws.Range("L112").Value = 14653.0005
$ws.Range("N112").Value = -16869.0005
$ws.Range("H129").Value = 4083.375
$ws.Range("I129").Value = 5173.75
$ws.Range("J129").Value = 2993
$ws.Range("K129").Value = 15521.25
$ws.Range("L129").Value = 8979
$ws.Range("M129").Value = -10521.25
$ws.Range("N129").Value = -18979
$ws.Range("H137").Value = 1304.3636
$ws.Range("I137").Value = 1349.8
$ws.Range("J137").Value = 850
$ws.Range("K137").Value = 4049.4
$ws.Range("L137").Value = 2550
$ws.Range("M137").Value = -1499.4
$ws.Range("N137").Value = -7650
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2201.25
$ws.Range("I102").Value = 2201.25
$ws.Range("K102").Value = 2201.25
$ws.Range("M102").Value = -579.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 84.166664
$ws.Range("I19").Value = 82.72727
$ws.Range("K19").Value = 82.72727
$ws.Range("M19").Value = 87.27273
$ws.Range("H21").Value = 515
$ws.Range("J21").Value = 515
$ws.Range("L21").Value = 515
$ws.Range("N21").Value = -985
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H24").Value = 84.166664
$ws.Range("I24").Value = 82.72727
$ws.Range("K24").Value = 82.72727
$ws.Range("M24").Value = 87.27273
$ws.Range("H31").Value = 2268.4
$ws.Range("I31").Value = 2030
$ws.Range("K31").Value = 2030
$ws.Range("M31").Value = -1735
$ws.Range("H34").Value = 2268.4
$ws.Range("I34").Value = 2030
$ws.Range("K34").Value = 2030
$ws.Range("M34").Value = -1828
$ws.Range("H35").Value = 1022.5
$ws.Range("I35").Value = 1022.5
$ws.Range("K35").Value = 1022.5
$ws.Range("M35").Value = -728.5
$ws.Range("H94").Value = 227290.8
$ws.Range("I94").Value = 283558
$ws.Range("J94").Value = 2222
$ws.Range("K94").Value = 283558
$ws.Range("L94").Value = 2222
$ws.Range("M94").Value = -283107
$ws.Range("N94").Value = -3124
$ws.Range("H141").Value = 35958.465
$ws.Range("J141").Value = 35958.465
$ws.Range("L141").Value = 35958.465
$ws.Range("N141").Value = -46318.465
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 134.1
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 165.125
$ws.Range("K2").Value = 60
$ws.Range("L2").Value = 990.75
$ws.Range("M2").Value = 53
$ws.Range("N2").Value = -1216.75
$ws.Range("H12").Value = 125.388885
$ws.Range("J12").Value = 59.5
$ws.Range("L12").Value = 178.5
$ws.Range("N12").Value = -524.5
$ws.Range("H22").Value = 2113.1052
$ws.Range("I22").Value = 1250
$ws.Range("J22").Value = 2214.647
$ws.Range("K22").Value = 3750
$ws.Range("L22").Value = 6643.941
$ws.Range("M22").Value = -3581
$ws.Range("N22").Value = -6981.941
$ws.Range("H26").Value = 400
$ws.Range("I26").Value = 400
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 1200
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -912
$ws.Range("N26").ClearContents()
$ws.Range("H27").Value = 2113.1052
$ws.Range("I27").Value = 1250
$ws.Range("J27").Value = 2214.647
$ws.Range("K27").Value = 3750
$ws.Range("L27").Value = 6643.941
$ws.Range("M27").Value = -3648
$ws.Range("N27").Value = -6847.941
$ws.Range("H55").Value = 5992.5
$ws.Range("J55").Value = 5992.5
$ws.Range("L55").Value = 17977.5
$ws.Range("N55").Value = -18331.5
$ws.Range("H59").Value = 2540.4
$ws.Range("I59").Value = 2005
$ws.Range("J59").Value = 2599.889
$ws.Range("K59").Value = 6015
$ws.Range("L59").Value = 7799.667
$ws.Range("M59").Value = -5475
$ws.Range("N59").Value = -8879.667000000001
$ws.Range("H76").Value = 15429.5
$ws.Range("J76").Value = 16363.637
$ws.Range("L76").Value = 49090.911
$ws.Range("N76").Value = -49856.911
$ws.Range("H79").Value = 15429.5
$ws.Range("J79").Value = 16363.637
$ws.Range("L79").Value = 49090.911
$ws.Range("N79").Value = -51742.911
$ws.Range("H81").Value = 2000
$ws.Range("I81").Value = 2000
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 6000
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -4877
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 2000
$ws.Range("I84").Value = 2000
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 18000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -12384
$ws.Range("N84").ClearContents()
$ws.Range("H98").Value = 1033.125
$ws.Range("I98").Value = 696
$ws.Range("K98").Value = 2088
$ws.Range("M98").Value = -590
$ws.Range("H109").Value = 3000
$ws.Range("I109").Value = 3000
$ws.Range("K109").Value = 9000
$ws.Range("M109").Value = -7960
$ws.Range("H113").Value = 1519.6
$ws.Range("J113").Value = 2099.6667
$ws.Range("L113").Value = 6299.000100000001
$ws.Range("N113").Value = -10639.0001
$ws.Range("H115").Value = 3239
$ws.Range("I115").Value = 2000
$ws.Range("J115").Value = 3486.8
$ws.Range("K115").Value = 6000
$ws.Range("L115").Value = 10460.4
$ws.Range("M115").Value = -4825
$ws.Range("N115").Value = -12810.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2099.75
$ws.Range("I80").Value = 2299.6667
$ws.Range("J80").Value = 1500
$ws.Range("K80").Value = 2299.6667
$ws.Range("L80").Value = 1500
$ws.Range("M80").Value = -1301.6667
$ws.Range("N80").Value = -3496
$ws.Range("H83").Value = 2099.75
$ws.Range("I83").Value = 2299.6667
$ws.Range("J83").Value = 1500
$ws.Range("K83").Value = 11498.3335
$ws.Range("L83").Value = 7500
$ws.Range("M83").Value = -6506.333500000001
$ws.Range("N83").Value = -17484
$ws.Range("H97").Value = 293
$ws.Range("I97").Value = 293
$ws.Range("K97").Value = 293
$ws.Range("M97").Value = 203
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H40").Value = 3895.0833
$ws.Range("I40").Value = 3325.3333
$ws.Range("J40").Value = 4464.8335
$ws.Range("K40").Value = 3325.3333
$ws.Range("L40").Value = 4464.8335
$ws.Range("M40").Value = -3189.3333
$ws.Range("N40").Value = -4736.8335
$ws.Range("H55").Value = 1050.2142
$ws.Range("I55").Value = 562.3333
$ws.Range("J55").Value = 1416.125
$ws.Range("K55").Value = 562.3333
$ws.Range("L55").Value = 1416.125
$ws.Range("M55").Value = -389.3333
$ws.Range("N55").Value = -1762.125
$ws.Range("H132").Value = 3452.3635
$ws.Range("I132").Value = 2747.25
$ws.Range("J132").Value = 5332.6665
$ws.Range("K132").Value = 8241.75
$ws.Range("L132").Value = 15997.9995
$ws.Range("M132").Value = -5711.75
$ws.Range("N132").Value = -21057.9995
$ws.Range("H136").Value = 3129.6667
$ws.Range("I136").Value = 3129.6667
$ws.Range("K136").Value = 9389.000100000001
$ws.Range("M136").Value = -6839.000100000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2035.4286
$ws.Range("I136").Value = 1424.5
$ws.Range("K136").Value = 4273.5
$ws.Range("M136").Value = -1723.5
